# The post "「結婚前、結婚後」" (row 422) was removed from posts.xlsx.
# Delete that entire row; Excel will automatically shift all subsequent
# rows (423-437) up by one and shrink the used range (dimension) from
# A1:C437 to A1:C436.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(422).Delete()
